# Added the Buy and sell project flow admin accept flow
# Appends 31 new "Solar Project" entries to column A of the
# "Project details" sheet (rows 68-98), continuing the existing list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project details")

$values = @(
    "Solar Project DA",
    "Solar Project CB",
    "Solar Project CC",
    "Solar Project CC",
    "Solar Project CC",
    "Solar Project BD",
    "Solar Project CC",
    "Solar Project DB",
    "Solar Project AD",
    "Solar Project CC",
    "Solar Project AD",
    "Solar Project CC",
    "Solar Project CA",
    "Solar Project DB",
    "Solar Project AA",
    "Solar Project DA",
    "Solar Project CA",
    "Solar Project AB",
    "Solar Project AC",
    "Solar Project AA",
    "Solar Project CA",
    "Solar Project BC",
    "Solar Project BD",
    "Solar Project CD",
    "Solar Project BD",
    "Solar Project BD",
    "Solar Project DD",
    "Solar Project BA",
    "Solar Project CC",
    "Solar Project CA",
    "Solar Project BB"
)

$startRow = 68
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
